$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New case rows appended after the existing data (21TRC08418 / Hemmeter).
# Each inner array is one row: columns A..K. "<<SKIP>>" leaves the cell
# untouched (blank/absent); "<<EMPTY>>" writes an explicit empty text value.
$newRows = @(
    @("21TRC08418", "Hemmeter", "DRIVING IN MARKED LANES", "4511.33", "MM", "Guilty", "Guilty", "`$ 15", "`$ 0", "None", "None"),  # row 291
    @("21TRC08418", "Hemmeter", "TURN AND STOP SIGNALS", "4511.39", "MM", "Dismissed", "<<SKIP>>", " ", " ", " ", " "),  # row 292
    @("21TRC08418", "Hemmeter", "OVI ALCOHOL / DRUGS 1ST", "4511.19A1A*", "M1", "Guilty", "Guilty", "`$ 250", "`$ 0", "180", "177"),  # row 293
    @("21TRC08418", "Hemmeter", "DRIVING IN MARKED LANES", "4511.33", "MM", "Guilty", "Guilty", "`$ 15", "`$ 0", "None", "None"),  # row 294
    @("21TRC08418", "Hemmeter", "TURN AND STOP SIGNALS", "4511.39", "MM", "Dismissed", "<<SKIP>>", " ", " ", " ", " "),  # row 295
    @("21TRC08418", "Hemmeter", "OVI ALCOHOL / DRUGS 1ST", "4511.19A1A*", "M1", "Guilty", "Guilty", "`$ 250", "`$ 0", "180", "177"),  # row 296
    @("21TRC08418", "Hemmeter", "DRIVING IN MARKED LANES", "4511.33", "MM", "Guilty", "Guilty", "`$ 15", "`$ 0", "None", "None"),  # row 297
    @("21TRC08418", "Hemmeter", "TURN AND STOP SIGNALS", "4511.39", "MM", "Dismissed", "<<SKIP>>", " ", " ", " ", " "),  # row 298
    @("21TRC08418", "Hemmeter", "OVI ALCOHOL / DRUGS 1ST", "4511.19A1A*", "M1", "Guilty", "Guilty", "`$ 100", "`$ 0", "180", "177"),  # row 299
    @("21TRC08418", "Hemmeter", "DRIVING IN MARKED LANES", "4511.33", "MM", "Guilty", "Guilty", "`$ 15", "`$ 0", "None", "None"),  # row 300
    @("21TRC08418", "Hemmeter", "TURN AND STOP SIGNALS", "4511.39", "MM", "Dismissed", "<<EMPTY>>", " ", " ", " ", " "),  # row 301
    @("21TRC08418", "Hemmeter", "OVI ALCOHOL / DRUGS 1ST", "4511.19A1A*", "M1", "Guilty", "Guilty", "`$ 100", "`$ 0", "180", "177"),  # row 302
)

$startRow = 291
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $newRows[$i]
    for ($col = 1; $col -le $rowVals.Count; $col++) {
        $val = $rowVals[$col - 1]
        if ($val -eq "<<SKIP>>") { continue }
        $c = $ws.Cells.Item($r, $col)
        $c.NumberFormat = "@"
        if ($val -eq "<<EMPTY>>") {
            $c.Value = "'"
        } else {
            $c.Value = $val
        }
    }
}